$d = $word.ActiveDocument

# Locate the paragraph that holds the stray "https://gitlab.com/jozefhajnala/jhaddins.git"
# hyperlink (a bullet item, numId 8) that sits right after the
# "– Just give me the package" bullet, just before the trailing empty
# paragraph / sectPr, and remove the whole paragraph (including its
# paragraph mark) so the surrounding structure collapses exactly like
# the diff shows.
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*gitlab.com/jozefhajnala/jhaddins.git*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
